$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cryptocurrency price/volume data refresh (GitHub Actions update)

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '64.947.85'
$ws.Range("E2").Value = '  -1.59%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.244.97'
$ws.Range("E3").Value = '  -0.57%  '

$ws.Range("E4").Value = '  +0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.07'
$ws.Range("E5").Value = '  +0.85%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '175.45'
$ws.Range("E6").Value = '  -1.48%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.629'
$ws.Range("E7").Value = '  +0.73%  '

$ws.Range("E8").Value = '  +0.08%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '3.242.53'
$ws.Range("E9").Value = '  -0.63%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.122'
$ws.Range("E10").Value = '  -2.25%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.76'
$ws.Range("E11").Value = '  +0.99%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.392'
$ws.Range("E12").Value = '  -2.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.822.09'
$ws.Range("E13").Value = '  -0.21%  '

$ws.Range("E14").Value = '  -2.98%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '65.119.76'
$ws.Range("E15").Value = '  -1.41%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '25.69'
$ws.Range("E16").Value = '  -2.73%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.238.17'
$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000158'
$ws.Range("E18").Value = '  -2.19%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '416.30'
$ws.Range("E19").Value = '  -3.84%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.37'
$ws.Range("E20").Value = '  -3.08%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.85'
$ws.Range("E21").Value = '  -2.10%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.18'
$ws.Range("E22").Value = '  -2.80%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.38%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '70.43'
$ws.Range("E24").Value = '  -1.80%  '

$ws.Range("E25").Value = '  -0.13%  '

$ws.Range("B26").Value = 'Kaspa'
$ws.Range("C26").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.202'
$ws.Range("E26").Value = '  +3.56%  '

$ws.Range("B27").Value = 'Polygon'
$ws.Range("C27").Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.497'
$ws.Range("E27").Value = '  -1.43%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0000110'
$ws.Range("E28").Value = '  -1.95%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.17'
$ws.Range("E29").Value = '  +3.52%  '

$ws.Range("E30").Value = '  +0.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.88'
$ws.Range("E31").Value = '  -2.50%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '21.86'
$ws.Range("E32").Value = '  -1.68%  '

$ws.Range("E33").Value = '  -0.04%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '5.05'
$ws.Range("E34").Value = '  -1.52%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.43'
$ws.Range("E35").Value = '  -2.36%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.15'
$ws.Range("E36").Value = '  -2.23%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '157.16'
$ws.Range("E37").Value = '  +0.06%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.40'
$ws.Range("E38").Value = '  -1.76%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.814.98'
$ws.Range("E39").Value = '  +1.82%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '1.74'
$ws.Range("E40").Value = '  -2.72%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '25.55'
$ws.Range("E41").Value = '  -3.40%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '4.22'
$ws.Range("E42").Value = '  -1.66%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.733'
$ws.Range("E43").Value = '  -5.30%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '39.40'
$ws.Range("E44").Value = '  -1.94%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.75'
$ws.Range("E45").Value = '  -4.35%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0626'
$ws.Range("E46").Value = '  -4.24%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '305.69'
$ws.Range("E47").Value = '  -4.51%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.17'
$ws.Range("E48").Value = '  -5.53%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '22.19'
$ws.Range("E49").Value = '  -4.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0262'
$ws.Range("E50").Value = '  -1.18%  '

$ws.Range("E51").Value = '  -1.06%  '
